$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J (year 2020) ---------------------------------------------

# Header cell J4: same look as I4 (bold header style with bottom border)
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 2020

# Data rows 5-10: copy the formatting from the matching column-I cell so the
# new column lines up visually with the rest of the table, then set values.
$ws.Range("I5").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("J5").Value = 370

$ws.Range("I6").Copy()
$ws.Range("J6").PasteSpecial(-4122)
$ws.Range("J6").Value = 5

$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = 5

$ws.Range("I8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Value = 20

$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = 19

$ws.Range("I10").Copy()
$ws.Range("J10").PasteSpecial(-4122)
$ws.Range("J10").Value = 73

# Row 3 is the thin separator row above the header; give the new J3 cell the
# same medium bottom border rule used across the rest of that row.
$ws.Range("J3").Borders.Item(9).Weight = -4138

# --- Updated figures in the existing 2019 column (I) -----------------------
$ws.Range("I8").Value = 42
$ws.Range("I9").Value = 30
$ws.Range("I10").Value = 62
